# Refresh cryptos list (Price / Volume(1h) columns) to the Wed Jul 31 05:31:08
# UTC 2024 GitHub Actions snapshot.
#
# Columns D (Price) and E (Volume(1h)) are stored as plain text in this
# sheet (t="inlineStr"), even when the price text happens to look like a
# number (e.g. "582.44", "0.999"). Excel's normal cell-entry parsing would
# silently convert a bare numeric-looking string typed into a General-
# formatted cell into a real number (and normalise "0.640" -> "0.64",
# "26.20" -> "26.2", etc.), which would NOT match the source data. To keep
# those values as literal text we prefix them with a leading apostrophe
# ('), exactly like typing a quote-prefixed entry in Excel - the apostrophe
# itself is not stored, it just forces text interpretation. Values that are
# already unambiguous as text (contain two dots, a trailing "%", stray
# spaces, etc.) are assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '65.734.87'
$ws.Range("E2").Value = '  -1.12%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '3.274.21'
$ws.Range("E3").Value = '  -0.90%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.05%  '

# Row 5: BNB
$ws.Range("D5").Value = '''582.44'
$ws.Range("E5").Value = '  +1.97%  '

# Row 6: Solana
$ws.Range("D6").Value = '''178.86'
$ws.Range("E6").Value = '  -1.73%  '

# Row 7: XRP
$ws.Range("D7").Value = '''0.640'
$ws.Range("E7").Value = '  +7.13%  '

# Row 8: USDC
$ws.Range("E8").Value = '  -0.06%  '

# Row 9: Dogecoin
$ws.Range("E9").Value = '  -3.61%  '

# Row 10: Toncoin
$ws.Range("D10").Value = '''6.73'
$ws.Range("E10").Value = '  +1.66%  '

# Row 11: Cardano
$ws.Range("E11").Value = '  -0.50%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '3.838.28'
$ws.Range("E12").Value = '  -1.09%  '

# Row 13: TRON
$ws.Range("E13").Value = '  -4.56%  '

# Row 14: WrappedBTC
$ws.Range("D14").Value = '65.772.72'
$ws.Range("E14").Value = '  -1.19%  '

# Row 15: Avalanche
$ws.Range("D15").Value = '''25.87'
$ws.Range("E15").Value = '  -4.56%  '

# Row 16: WrappedEther
$ws.Range("D16").Value = '3.278.75'
$ws.Range("E16").Value = '  -0.96%  '

# Row 17: ShibaInu
$ws.Range("D17").Value = '''0.0000162'
$ws.Range("E17").Value = '  -3.14%  '

# Row 18: BitcoinCash
$ws.Range("D18").Value = '''424.47'
$ws.Range("E18").Value = '  -1.20%  '

# Row 19: Chainlink
$ws.Range("D19").Value = '''13.14'
$ws.Range("E19").Value = '  -3.74%  '

# Row 20: Polkadot
$ws.Range("D20").Value = '''5.47'
$ws.Range("E20").Value = '  -3.40%  '

# Row 21: Uniswap
$ws.Range("D21").Value = '''7.34'
$ws.Range("E21").Value = '  -3.56%  '

# Row 22: Dai
$ws.Range("E22").Value = '  +0.01%  '

# Row 23: Litecoin
$ws.Range("E23").Value = '  -2.87%  '

# Row 24: LEO
$ws.Range("D24").Value = '''5.67'
$ws.Range("E24").Value = '  +0.20%  '

# Row 25: WrappedeETH
$ws.Range("D25").Value = '3.419.80'
$ws.Range("E25").Value = '  -0.96%  '

# Row 26: Polygon
$ws.Range("D26").Value = '''0.507'
$ws.Range("E26").Value = '  -1.50%  '

# Row 27: Kaspa
$ws.Range("E27").Value = '  +1.41%  '

# Row 28: PEPE
$ws.Range("D28").Value = '''0.0000112'
$ws.Range("E28").Value = '  -5.13%  '

# Row 29: InternetComputer(DFINITY)
$ws.Range("D29").Value = '''8.82'
$ws.Range("E29").Value = '  -2.03%  '

# Row 30: Binance-PegBSC-USD
$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  -0.16%  '

# Row 31: PancakeSwap
$ws.Range("E31").Value = '  -0.61%  '

# Row 32: EthereumClassic
$ws.Range("D32").Value = '''22.16'
$ws.Range("E32").Value = '  -2.55%  '

# Row 33: USDe
$ws.Range("E33").Value = '  +0.08%  '

# Row 34: NEARProtocol
$ws.Range("D34").Value = '''5.12'
$ws.Range("E34").Value = '  -3.61%  '

# Row 35: Aptos
$ws.Range("D35").Value = '''6.55'
$ws.Range("E35").Value = '  -3.19%  '

# Row 36: Fetch.AI
$ws.Range("E36").Value = '  -4.19%  '

# Row 37: Monero
$ws.Range("D37").Value = '''158.93'
$ws.Range("E37").Value = '  -0.61%  '

# Row 38: ImmutableX
$ws.Range("D38").Value = '''1.42'
$ws.Range("E38").Value = '  -5.44%  '

# Row 39: Stacks
$ws.Range("D39").Value = '''1.79'
$ws.Range("E39").Value = '  -3.03%  '

# Row 40: EnergySwap
$ws.Range("D40").Value = '''26.20'
$ws.Range("E40").Value = '  -3.34%  '

# Row 41: Maker
$ws.Range("D41").Value = '2.784.21'
$ws.Range("E41").Value = '  -0.28%  '

# Row 42: Mantle
$ws.Range("D42").Value = '''0.764'
$ws.Range("E42").Value = '  -3.05%  '

# Row 43: Filecoin
$ws.Range("D43").Value = '''4.30'
$ws.Range("E43").Value = '  -2.97%  '

# Row 44: OKB
$ws.Range("D44").Value = '''39.90'
$ws.Range("E44").Value = '  -0.64%  '

# Row 45: Hedera
$ws.Range("D45").Value = '''0.0654'
$ws.Range("E45").Value = '  -2.85%  '

# Row 46: RenderToken
$ws.Range("D46").Value = '''5.84'
$ws.Range("E46").Value = '  -5.28%  '

# Row 47: dogwifhat
$ws.Range("D47").Value = '''2.27'
$ws.Range("E47").Value = '  -3.15%  '

# Row 48: Bittensor
$ws.Range("D48").Value = '''313.72'
$ws.Range("E48").Value = '  -2.02%  '

# Row 49: InjectiveProtocol
$ws.Range("D49").Value = '''22.99'
$ws.Range("E49").Value = '  -5.44%  '

# Row 50: VeChain
$ws.Range("D50").Value = '''0.0266'
$ws.Range("E50").Value = '  -2.01%  '

# Row 51: Stellar
$ws.Range("E51").Value = '  +3.18%  '
